$d = $word.ActiveDocument

# Replace the "ABRIL / MAIO" text with "CAMPO MES" (split across two runs
# in the underlying OOXML, but from the Find/Replace perspective it is a
# straightforward text substitution).
$d.Content.Find.Execute("ABRIL / MAIO", $true, $false, $false, $false, $false,
                         $true, 1, $false, "CAMPO MES", 2)
